$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.05052846209076733
$ws.Range("D2").Value = 0.008943546438445082
$ws.Range("E2").Value = 0.4304237066138796
$ws.Range("F2").Value = 1.038957906139956
$ws.Range("G2").Value = 0.002342237695844014
$ws.Range("O2").Value = 3.581084263747925

$ws.Range("C3").Value = 0.04476206481498934
$ws.Range("D3").Value = 0.007995848526700655
$ws.Range("E3").Value = 0.3748685171765942
$ws.Range("F3").Value = 0.9555018548654175
$ws.Range("G3").Value = 0.002348116949748624
$ws.Range("O3").Value = 3.290152896055872

$ws.Range("C4").Value = 0.04124261813029761
$ws.Range("D4").Value = 0.00742011289122857
$ws.Range("E4").Value = 0.3409164150049975
$ws.Range("F4").Value = 0.904935892906451
$ws.Range("G4").Value = 0.002351912108086821
$ws.Range("O4").Value = 3.113860341885584

$ws.Range("C5").Value = 0.03981353979432356
$ws.Range("D5").Value = 0.007186965552236302
$ws.Range("E5").Value = 0.3271169890139163
$ws.Range("F5").Value = 0.8844966276862323
$ws.Range("G5").Value = 0.002353505441242354
$ws.Range("O5").Value = 3.042596276585186

$ws.Range("C6").Value = 0.0395765463820652
$ws.Range("D6").Value = 0.007148338117023201
$ws.Range("E6").Value = 0.3248276933735923
$ws.Range("F6").Value = 0.8811126802678615
$ws.Range("G6").Value = 0.002353772843926084
$ws.Range("O6").Value = 3.030797410986111

$ws.Range("C7").Value = 0.04122332458312883
$ws.Range("D7").Value = 0.007416962736066068
$ws.Range("E7").Value = 0.340730168809614
$ws.Range("F7").Value = 0.9046595708997387
$ws.Range("G7").Value = 0.002351933406695056
$ws.Range("O7").Value = 3.112896930884801

$ws.Range("C8").Value = 0.04853570168536692
$ws.Range("D8").Value = 0.008615455172950703
$ws.Range("E8").Value = 0.4112330468379781
$ws.Range("F8").Value = 1.010039843139111
$ws.Range("G8").Value = 0.002344226522264665
$ws.Range("O8").Value = 3.480278027715428

$ws.Range("C9").Value = 0.06305247523306434
$ws.Range("D9").Value = 0.01101820639376427
$ws.Range("E9").Value = 0.5509232592974627
$ws.Range("F9").Value = 1.222218960642778
$ws.Range("G9").Value = 0.002330574835433471
$ws.Range("O9").Value = 4.219864929822563

$ws.Range("C10").Value = 0.07384021095187165
$ws.Range("D10").Value = 0.01282084040147424
$ws.Range("E10").Value = 0.6546786155740421
$ws.Range("F10").Value = 1.381719871913816
$ws.Range("G10").Value = 0.00232142388172245
$ws.Range("O10").Value = 4.775787231787604

$ws.Range("C11").Value = 0.07877756143632553
$ws.Range("D11").Value = 0.01365020014087293
$ws.Range("E11").Value = 0.7021797081819585
$ws.Range("F11").Value = 1.455119823127035
$ws.Range("G11").Value = 0.002317449146469602
$ws.Range("O11").Value = 5.031610505972708

$ws.Range("C12").Value = 0.0806517566815188
$ws.Range("D12").Value = 0.0139656954419678
$ws.Range("E12").Value = 0.7202150718743354
$ws.Range("F12").Value = 1.483039761349005
$ws.Range("G12").Value = 0.002315970862932314
$ws.Range("O12").Value = 5.12892072231125

$ws.Range("C13").Value = 0.08024791027854405
$ws.Range("D13").Value = 0.01389768257124757
$ws.Range("E13").Value = 0.7163286418724084
$ws.Range("F13").Value = 1.477021082648264
$ws.Range("G13").Value = 0.002316288045592597
$ws.Range("O13").Value = 5.107943632011938

$ws.Range("C14").Value = 0.07893166058768486
$ws.Range("D14").Value = 0.01367612683552721
$ws.Range("E14").Value = 0.7036625070413862
$ws.Range("F14").Value = 1.45741428624936
$ws.Range("G14").Value = 0.002317326989741375
$ws.Range("O14").Value = 5.039607467092196

$ws.Range("C15").Value = 0.07812601590281076
$ws.Range("D15").Value = 0.01354060718945505
$ws.Range("E15").Value = 0.6959104778898535
$ws.Range("F15").Value = 1.445420948117032
$ws.Range("G15").Value = 0.002317966865919552
$ws.Range("O15").Value = 4.997806726364843

$ws.Range("C16").Value = 0.07351816499674158
$ws.Range("D16").Value = 0.01276683546616653
$ws.Range("E16").Value = 0.6515807523204842
$ws.Range("F16").Value = 1.376940287943881
$ws.Range("G16").Value = 0.002321687408983693
$ws.Range("O16").Value = 4.759128753694711

$ws.Range("C17").Value = 0.07069922541678864
$ws.Range("D17").Value = 0.01229460716505315
$ws.Range("E17").Value = 0.6244663244896742
$ws.Range("F17").Value = 1.335148049666202
$ws.Range("G17").Value = 0.002324017883126265
$ws.Range("O17").Value = 4.613468198980513

$ws.Range("C18").Value = 0.06908065205573166
$ws.Range("D18").Value = 0.01202386667504385
$ws.Range("E18").Value = 0.6088990123417233
$ws.Range("F18").Value = 1.311189294299879
$ws.Range("G18").Value = 0.002325376024032553
$ws.Range("O18").Value = 4.52996316860856

$ws.Range("C19").Value = 0.06853310624022413
$ws.Range("D19").Value = 0.01193234569974777
$ws.Range("E19").Value = 0.6036329010448469
$ws.Range("F19").Value = 1.303090729238818
$ws.Range("G19").Value = 0.002325838914997291
$ws.Range("O19").Value = 4.50173662271493

$ws.Range("C20").Value = 0.07099901434227718
$ws.Range("D20").Value = 0.01234478568321151
$ws.Range("E20").Value = 0.6273497534900372
$ws.Range("F20").Value = 1.339588697244295
$ws.Range("G20").Value = 0.002323767968160808
$ws.Range("O20").Value = 4.628945444994827

$ws.Range("C21").Value = 0.07931815036987189
$ws.Range("D21").Value = 0.01374116343710341
$ws.Range("E21").Value = 0.7073815280193969
$ws.Range("F21").Value = 1.463169854790664
$ws.Range("G21").Value = 0.002317021099142336
$ws.Range("O21").Value = 5.059667523342227

$ws.Range("C22").Value = 0.08478168316915458
$ws.Range("D22").Value = 0.014662184576828
$ws.Range("E22").Value = 0.7599671772600232
$ws.Range("F22").Value = 1.544667555240096
$ws.Range("G22").Value = 0.002312768120793171
$ws.Range("O22").Value = 5.343714382450344

$ws.Range("C23").Value = 0.08186319461671587
$ws.Range("D23").Value = 0.01416981774027448
$ws.Range("E23").Value = 0.7318741492869805
$ws.Range("F23").Value = 1.501102581173711
$ws.Range("G23").Value = 0.00231502375701132
$ws.Range("O23").Value = 5.191875664919166

$ws.Range("C24").Value = 0.07086347332972309
$ws.Range("D24").Value = 0.0123220976466385
$ws.Range("E24").Value = 0.6260460897340892
$ws.Range("F24").Value = 1.337580868464897
$ws.Range("G24").Value = 0.002323880897583677
$ws.Range("O24").Value = 4.621947444414047

$ws.Range("C25").Value = 0.0591048735734887
$ws.Range("D25").Value = 0.01036203391903001
$ws.Range("E25").Value = 0.5129541931938206
$ws.Range("F25").Value = 1.164202196230377
$ws.Range("G25").Value = 0.002334112758844178
$ws.Range("O25").Value = 4.017647422946709

